$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that gets bumped by
# one day (automatic update) for every data row (rows 2 through 23).
for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 3).Value = 46082
}
